$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Formula = "=(0.0438 + 0.0412 + 0.0403) / 3 * 100"
$ws.Range("C20").Formula = "=(0.0365 + 0.0322 + 0.0132) / 3 * 100"
$ws.Range("E21").Formula = "=0.0998 * 100"
$ws.Range("E22").Formula = "=0.1108*100"

$ws.Range("J16").Select()
